$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 824, shifting existing rows 824+ down by one
$ws.Rows(824).Insert()

# Populate the newly inserted row 824 with the new data point
# Force column A to be stored as text (not auto-converted to a date serial)
$ws.Cells.Item(824, 1).NumberFormat = "@"
$ws.Cells.Item(824, 1).Value = "2026/02/15"
$ws.Cells.Item(824, 2).Value = "日"
$ws.Cells.Item(824, 3).Value = 13
$ws.Cells.Item(824, 4).Value = 201
